$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workflow")

# Insert a new row at position 7 (pushes the old rows 7 "random effects
# intercepts" and 8 "random slopes" down to rows 8 and 9).
$ws.Rows.Item(7).Insert() | Out-Null

# Row 4 ("Find best seasonality spec?") now defaults to FALSE and the
# R_name changes from pick_fft to search_seasonality.
$ws.Cells.Item(4, 2).Value = $false
$ws.Cells.Item(4, 3).Value = "search_seasonality"

# Newly inserted row 7: no description, defaults to FALSE, R_name
# search_randoms.
$ws.Cells.Item(7, 2).Value = $false
$ws.Cells.Item(7, 3).Value = "search_randoms"

# Move the active selection to the newly inserted cell.
$ws.Range("B7").Select() | Out-Null
